{"js": "// The edit removes four paragraphs that used to sit right after the\n// \"LOT2004: Bioqu\u00edmica (Requisito fraco)\" requirement line:\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) another empty paragraph\n//   4) an empty paragraph carrying a page-break-before\n// The paragraph right after that block (a plain empty paragraph) is kept,\n// as is everything before the anchor line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOT2004\" requirement paragraph \u2014 the anchor right before the\n// block that needs to be removed.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOT2004\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // The text of the paragraphs to drop, in document order, starting right\n  // after the anchor paragraph.\n  const removableTexts = [\"\", \"Ver no Jupiter Salvar em pdf Salvar em docx\", \"\", \"\"];\n\n  const toDelete = [];\n  let cursor = anchorIndex + 1;\n  for (let k = 0; k < removableTexts.length && cursor < items.length; k++, cursor++) {\n    if (items[cursor].text !== removableTexts[k]) {\n      // Structure didn't match what we expected \u2014 bail out instead of\n      // deleting the wrong content.\n      toDelete.length = 0;\n      break;\n    }\n    toDelete.push(items[cursor]);\n  }\n\n  // Delete from the end backwards so earlier indices stay valid.\n  for (let k = toDelete.length - 1; k >= 0; k--) {\n    toDelete[k].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# The edit removes four paragraphs that used to sit right after the\n# \"LOT2004: Bioqu\u00edmica (Requisito fraco)\" requirement line:\n#   1) an empty paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) another empty paragraph\n#   4) an empty paragraph carrying a page-break-before\n# The paragraph right after that block (a plain empty paragraph) is kept,\n# as is everything before the anchor line.\n\n$d = $word.ActiveDocument\n\n# Locate the \"LOT2004\" requirement paragraph \u2014 the anchor right before the\n# block that needs to be removed.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*LOT2004*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 1) {\n    # The text of the paragraphs to drop, in document order, starting right\n    # after the anchor paragraph.\n    $removableTexts = @(\"\", \"Ver no Jupiter Salvar em pdf Salvar em docx\", \"\", \"\")\n\n    $matchCount = 0\n    for ($k = 0; $k -lt $removableTexts.Length; $k++) {\n        $idx = $anchorIndex + 1 + $k\n        if ($idx -gt $d.Paragraphs.Count) { break }\n        $text = $d.Paragraphs.Item($idx).Range.Text.TrimEnd([char]13, [char]7)\n        if ($text -ne $removableTexts[$k]) { break }\n        $matchCount = $matchCount + 1\n    }\n\n    if ($matchCount -eq $removableTexts.Length) {\n        $startPara = $d.Paragraphs.Item($anchorIndex + 1)\n        $endPara = $d.Paragraphs.Item($anchorIndex + $matchCount)\n        $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n        $delRange.Delete()\n    }\n}\n"}
